$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "happy"
$ws.Range("B9").Value = "vui mừng"
